# "round 1 final entries"
#
# A new player, Jordan Hansen, is inserted into the alphabetically-sorted
# leaderboard (column B) at row 22. That pushes every row from the old
# row 22 through the old last row (45) down by one, so the table grows
# from 44 entries (rows 2-45) to 45 entries (rows 2-46). The brand new
# row (22) only gets a rank number + name - no Round-1/Total values yet,
# and no bold-name style applied (it hasn't been formatted like the rest
# of the table yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 is brand new - give it the same formatting as the row above it
# before filling in values (existing rows 23-45 already carry the right
# per-column formatting since we are only overwriting their content).
$ws.Range("A45:D45").Copy($ws.Range("A46:D46"))

# Shift rows 22-45 down into 23-46, moving the actual cell content (not
# retyping it) so existing shared-string references stay intact. Walk
# bottom-up so a row's old content is read before it gets overwritten.
for ($r = 45; $r -ge 22; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value2
}
$ws.Cells.Item(46, 1).Value = 45

# The freshly inserted row: only rank + name, default (non-bold) style,
# and no Round 1 / Total cells at all yet.
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 2).Value = "Jordan Hansen"
$ws.Cells.Item(22, 3).ClearContents()
$ws.Cells.Item(22, 4).ClearContents()

# View state: sheet scrolled down a bit, selection now spans the extra row.
$ws.Range("A2:A46").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A2").Activate()

# Turn on the AutoFilter over the now-46-row table and record the
# sort-by-name state, mirroring what Excel persists after sorting /
# filtering this range.
$ws.Range("A1:D46").AutoFilter() | Out-Null

$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$46")
$fd.Visible = $false
